$d = $word.ActiveDocument

# The '_GoBack' bookmark marks the point right after "...is clicked. " at
# the end of the document (a zero-length bookmark). We'll rebuild the
# content from that point forward: split off a new empty paragraph, then
# a paragraph with the "Within..." text (with the _GoBack bookmark
# re-inserted partway through, matching where the user's cursor ended up),
# then a final paragraph with "The for the dots..." text.

$bm = $d.Bookmarks.Item("_GoBack")
$pos = $bm.Start
$bm.Delete()

# End the current paragraph (creates a new, empty paragraph after it).
$r = $d.Range($pos, $pos)
$r.InsertParagraphAfter()
$pos = $pos + 1

# End that new paragraph too, so it stays empty, and start the next one.
$r = $d.Range($pos, $pos)
$r.InsertParagraphAfter()
$pos = $pos + 1

# Paragraph: "With" + "in the HTML code, is calls the JavaS" + bookmark +
# "cript function 'plus slides' so that when the arrow is clicked it will
# go to the appropriate slide. "
$r = $d.Range($pos, $pos)
$r.InsertAfter("With")
$pos = $r.End

$r = $d.Range($pos, $pos)
$r.InsertAfter("in the HTML code, is calls the JavaS")
$pos = $r.End

$d.Bookmarks.Add("_GoBack", $d.Range($pos, $pos)) | Out-Null

$r = $d.Range($pos, $pos)
$r.InsertAfter("cript function ‘plus slides’ so that when the arrow is clicked it will go to the appropriate slide. ")
$pos = $r.End

$r = $d.Range($pos, $pos)
$r.InsertParagraphAfter()
$pos = $pos + 1

# Final paragraph: "The for the dots the function 'current slide' is
# called so that depending on what 'dot' is clicked, the slide displayed
# will mirror that."
$r = $d.Range($pos, $pos)
$r.InsertAfter("The for the dots the function ‘current slide’ is called so that depending on what ‘dot’ is clicked, the slide displayed will mirror that.")
